$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2699.7856
$ws.Range("I70").Value = 1756
$ws.Range("J70").Value = 2857.0833
$ws.Range("K70").Value = 5268
$ws.Range("L70").Value = 8571.249899999999
$ws.Range("M70").Value = -4998
$ws.Range("N70").Value = -9111.249899999999

$ws.Range("H73").Value = 2699.7856
$ws.Range("I73").Value = 1756
$ws.Range("J73").Value = 2857.0833
$ws.Range("K73").Value = 5268
$ws.Range("L73").Value = 8571.249899999999
$ws.Range("M73").Value = -4332
$ws.Range("N73").Value = -10443.2499

$ws.Range("H97").Value = 3848.6
$ws.Range("J97").Value = 3848.6
$ws.Range("L97").Value = 11545.8
$ws.Range("N97").Value = -12537.8

$ws.Range("H107").Value = 987.7646999999999
$ws.Range("I107").Value = 723
$ws.Range("J107").Value = 1223.1111
$ws.Range("K107").Value = 723
$ws.Range("L107").Value = 1223.1111
$ws.Range("M107").Value = 1197
$ws.Range("N107").Value = -5063.1111

$ws.Range("H135").Value = 1380.4667
$ws.Range("J135").Value = 1041
$ws.Range("L135").Value = 9369
$ws.Range("N135").Value = -14439

$ws.Range("H138").Value = 3413.16
$ws.Range("J138").Value = 3983.2632
$ws.Range("L138").Value = 11949.7896
$ws.Range("N138").Value = -22229.7896

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2320.1
$ws.Range("I74").Value = 2025.7778
$ws.Range("J74").Value = 4969
$ws.Range("K74").Value = 2025.7778
$ws.Range("L74").Value = 4969
$ws.Range("M74").Value = -1151.7778
$ws.Range("N74").Value = -6717

$ws.Range("H77").Value = 2320.1
$ws.Range("I77").Value = 2025.7778
$ws.Range("J77").Value = 4969
$ws.Range("K77").Value = 10128.889
$ws.Range("L77").Value = 24845
$ws.Range("M77").Value = -5760.889000000001
$ws.Range("N77").Value = -33581

$ws.Range("H110").Value = 3064.1428
$ws.Range("I110").Value = 2089.8
$ws.Range("K110").Value = 2089.8
$ws.Range("M110").Value = -44.80000000000018

$ws.Range("H132").Value = 2171
$ws.Range("I132").Value = 1984.4667
$ws.Range("K132").Value = 5953.4001
$ws.Range("M132").Value = -3423.4001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 5433.3335
$ws.Range("I96").Value = 5433.3335
$ws.Range("K96").Value = 5433.3335
$ws.Range("M96").Value = -2687.3335

$ws.Range("H107").Value = 13290.5
$ws.Range("I107").Value = 5322.8887
$ws.Range("J107").Value = 84999
$ws.Range("K107").Value = 5322.8887
$ws.Range("L107").Value = 84999
$ws.Range("M107").Value = -3402.8887
$ws.Range("N107").Value = -88839

$ws.Range("H134").Value = 6318
$ws.Range("I134").Value = 2049.3333
$ws.Range("J134").Value = 14001.6
$ws.Range("K134").Value = 6147.999899999999
$ws.Range("L134").Value = 42004.8
$ws.Range("M134").Value = -3612.999899999999
$ws.Range("N134").Value = -47074.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 502.25
$ws.Range("I16").Value = 502.25
$ws.Range("K16").Value = 502.25
$ws.Range("M16").Value = -215.25

$ws.Range("H31").Value = 6193.2856
$ws.Range("I31").Value = 6870.8
$ws.Range("K31").Value = 6870.8
$ws.Range("M31").Value = -6575.8

$ws.Range("H34").Value = 6193.2856
$ws.Range("I34").Value = 6870.8
$ws.Range("K34").Value = 6870.8
$ws.Range("M34").Value = -6668.8

$ws.Range("H99").Value = 866.6667
$ws.Range("J99").Value = 800
$ws.Range("L99").Value = 800
$ws.Range("N99").Value = -3796

$ws.Range("H105").Value = 538.5
$ws.Range("I105").Value = 533
$ws.Range("K105").Value = 533
$ws.Range("M105").Value = 1214

$ws.Range("H107").Value = 13000
$ws.Range("J107").Value = 13000
$ws.Range("L107").Value = 13000
$ws.Range("N107").Value = -16840

$ws.Range("H113").Value = 502.25
$ws.Range("I113").Value = 502.25
$ws.Range("K113").Value = 502.25
$ws.Range("M113").Value = 1667.75

$ws.Range("H126").Value = 866.6667
$ws.Range("J126").Value = 800
$ws.Range("L126").Value = 2400
$ws.Range("N126").Value = -7340

$ws.Range("H132").Value = 2105
$ws.Range("I132").Value = 1706
$ws.Range("K132").Value = 5118
$ws.Range("M132").Value = -2588

$ws.Range("H134").Value = 2317.6
$ws.Range("I134").Value = 1897
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 5691
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -3156
$ws.Range("N134").Value = -17070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 999.5
$ws.Range("I3").Value = 999.5
$ws.Range("K3").Value = 2998.5
$ws.Range("M3").Value = -2886.5

$ws.Range("H5").Value = 1806.7142
$ws.Range("I5").Value = 1120.5
$ws.Range("J5").Value = 1968.1765
$ws.Range("K5").Value = 3361.5
$ws.Range("L5").Value = 5904.529500000001
$ws.Range("M5").Value = -3249.5
$ws.Range("N5").Value = -6128.529500000001

$ws.Range("H34").Value = 1120
$ws.Range("I34").Value = 630.25
$ws.Range("K34").Value = 1890.75
$ws.Range("M34").Value = -1806.75

$ws.Range("H135").Value = 1806.7142
$ws.Range("I135").Value = 1120.5
$ws.Range("J135").Value = 1968.1765
$ws.Range("K135").Value = 10084.5
$ws.Range("L135").Value = 17713.5885
$ws.Range("M135").Value = -7549.5
$ws.Range("N135").Value = -22783.5885

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 2510000
$ws.Range("J52").Value = 2510000
$ws.Range("L52").Value = 2510000
$ws.Range("N52").Value = -2510518

$ws.Range("H113").Value = 925
$ws.Range("I113").Value = 750
$ws.Range("J113").Value = 1100
$ws.Range("K113").Value = 750
$ws.Range("L113").Value = 1100
$ws.Range("M113").Value = 1420
$ws.Range("N113").Value = -5440

$ws.Range("H132").Value = 4170
$ws.Range("I132").Value = 3255.5
$ws.Range("K132").Value = 9766.5
$ws.Range("M132").Value = -7236.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3016.5
$ws.Range("I61").Value = 2524.75
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 2524.75
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -2322.75
$ws.Range("N61").Value = -4404

$ws.Range("H113").Value = 3016.5
$ws.Range("I113").Value = 2524.75
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 2524.75
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -354.75
$ws.Range("N113").Value = -8340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 17001.5
$ws.Range("I3").Value = 14003
$ws.Range("J3").Value = 20000
$ws.Range("K3").Value = 14003
$ws.Range("L3").Value = 20000
$ws.Range("M3").Value = -13889
$ws.Range("N3").Value = -20228

$ws.Range("H11").Value = 30004
$ws.Range("I11").Value = 30004
$ws.Range("K11").Value = 30004
$ws.Range("M11").Value = -29862

$ws.Range("H81").Value = 1296.2222
$ws.Range("I81").Value = 833.25
$ws.Range("J81").Value = 5000
$ws.Range("K81").Value = 1666.5
$ws.Range("L81").Value = 10000
$ws.Range("M81").Value = -605.5
$ws.Range("N81").Value = -12122

$ws.Range("H84").Value = 1296.2222
$ws.Range("I84").Value = 833.25
$ws.Range("J84").Value = 5000
$ws.Range("K84").Value = 8332.5
$ws.Range("L84").Value = 50000
$ws.Range("M84").Value = -3028.5
$ws.Range("N84").Value = -60608
